$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Unmerge the cells currently sitting on rows 8 & 9 (totals + footer) ---
#        before relocating them, otherwise their non-anchor cells become stuck.
$ws.Range("P8:Q8").UnMerge()
$ws.Range("A9:F9").UnMerge()
$ws.Range("G9:I9").UnMerge()
$ws.Range("K9:Q9").UnMerge()

# --- 2. Relocate existing rows 8 (totals) and 9 (footer) down to 10 and 11 ---
#        bottom-up so the moves don't clobber each other.
$ws.Range("A9:Q9").Cut($ws.Range("A11:Q11"))
$ws.Range("A8:Q8").Cut($ws.Range("A10:Q10"))

# --- 3. Populate new rows 8 and 9 with formatting copied from row 7 ---
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)
$ws.Range("A7:Q7").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)

# Row heights for the two new product rows.
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5

# --- 4. Row 8 values: E-MOX 500MG 16 CAPS ---
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "E-MOX 500MG 16 CAPS"
$ws.Range("H8").Value = "0:1"
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = "1"
$ws.Range("L8").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"
$ws.Range("N8").Value = "40.00"
$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "20.0000"
$ws.Range("P8").NumberFormat = "0.00"
$ws.Range("Q8").Value = "0:1"

# --- 5. Row 9 values: VONDALOUS 20MG 20 F.C. TAB ---
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "VONDALOUS 20MG 20 F.C. TAB"
$ws.Range("H9").Value = "0:1"
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1"
$ws.Range("L9").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"
$ws.Range("N9").Value = "138.00"
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "69.0000"
$ws.Range("P9").NumberFormat = "0.00"
$ws.Range("Q9").Value = "0:1"

# --- 6. Update totals row (now row 10): 36 -> 125 ---
$ws.Range("P10").Value = 125

# --- 7. Refresh the generated-at timestamp in the footer row (now row 11) ---
$ws.Range("A11").Value = "Sunday, 3 August, 2025 9:53 AM"

# --- 8. Re-merge the cells for the two new rows and the relocated rows ---
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()
$ws.Range("P10:Q10").Merge()
$ws.Range("A11:F11").Merge()
$ws.Range("G11:I11").Merge()
$ws.Range("K11:Q11").Merge()
